$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.131.77'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.66%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.268.86'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.72%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '178.80'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.623'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.61%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -3.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.71'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.401'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.831.27'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.59%  '
$ws.Range('E13').Value = '  -3.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '66.167.31'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.46'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.55%  '
$ws.Range('E16').Value = '  -3.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.267.50'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '432.76'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.72%  '
$ws.Range('E19').Value = '  -2.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.96'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.99%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.411.26'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.505'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.40%  '
$ws.Range('E26').Value = '  +3.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000113'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.89'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('E30').Value = '  -2.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.25'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.14'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.79%  '
$ws.Range('E34').Value = '  -3.79%  '
$ws.Range('E35').Value = '  -4.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '156.78'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.66%  '
$ws.Range('E37').Value = '  -6.35%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.79'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.31%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '26.44'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.773.65'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.773'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.30'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.21'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.01'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0655'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.31'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.90%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '321.96'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.30'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0266'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.69%  '
$ws.Range('E50').Value = '  +1.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.999'
$ws.Range('D51').Style = 'Normal'
